# Generate Report for Handoff
# This swaps the "320a..." and "524b..." rows on all three sheets
# (the 524b item finished handback earlier, so it now sorts into row 2;
#  the 320a item has moved on to a new handoff cycle, so it now sits in
#  row 3 with fresh status/dates), and records the new handoff metadata
#  for the 320a item.

$wb = $excel.ActiveWorkbook

$urlTest0_320a = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d52cd53a6302f28979f4afec74a3120ec385c129/e2e/320a0865-d16b-4ab5-9df0-33233ac5f952.md"
$urlTest0_524b = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d52cd53a6302f28979f4afec74a3120ec385c129/e2e/524b8130-6d89-47c7-a944-7f48f88e90db.md"
$urlZhcn_320a  = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/50fb2878f6d1c7d2283296fe32a139546e551b81/e2e/320a0865-d16b-4ab5-9df0-33233ac5f952.md"
$urlZhcn_524b  = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/50fb2878f6d1c7d2283296fe32a139546e551b81/e2e/524b8130-6d89-47c7-a944-7f48f88e90db.md"
$urlDede_320a  = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/25d6cba1ccb363225b1269e7e7d833a73b667b29/e2e/320a0865-d16b-4ab5-9df0-33233ac5f952.md"
$urlDede_524b  = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/25d6cba1ccb363225b1269e7e7d833a73b667b29/e2e/524b8130-6d89-47c7-a944-7f48f88e90db.md"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d52cd53a6302f28979f4afec74a3120ec385c129/e2e/320a0865-d16b-4ab5-9df0-33233ac5f952.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a73dac993c2a5793e1ac0ee11d5fc12f5a169c18/e2e/320a0865-d16b-4ab5-9df0-33233ac5f952.md."

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "524b8130-6d89-47c7-a944-7f48f88e90db.md"
$ws1.Range("A3").Value = "320a0865-d16b-4ab5-9df0-33233ac5f952.md"

$ws1.Range("E3").Value = "Ready for handoff"
$ws1.Range("F3").Value = "Ready for handoff"
$ws1.Range("G3").Value = "2016-10-13 13:19:24"

# Hyperlinks: rebuild so the ref/target assignment (rId2 -> 320a url,
# rId3 -> 524b url) is preserved but the displayed text is swapped.
$ws1.Range("B2").Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), $urlTest0_320a, [Type]::Missing, [Type]::Missing, "e2e\524b8130-6d89-47c7-a944-7f48f88e90db.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("B3"), $urlTest0_524b, [Type]::Missing, [Type]::Missing, "e2e\320a0865-d16b-4ab5-9df0-33233ac5f952.md") | Out-Null

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "524b8130-6d89-47c7-a944-7f48f88e90db.md"
$ws2.Range("I2").Value = "524b8130-6d89-47c7-a944-7f48f88e90db.md"
$ws2.Range("G2").Value = "524b8130-6d89-47c7-a944-7f48f88e90db.29f77a1baee9a4f82de0e54920eab7c2a3320007.zh-cn.xlf"
$ws2.Range("J2").Value = "524b8130-6d89-47c7-a944-7f48f88e90db.29f77a1baee9a4f82de0e54920eab7c2a3320007.zh-cn.xlf"

$ws2.Range("A3").Value = "320a0865-d16b-4ab5-9df0-33233ac5f952.md"
$ws2.Range("I3").Value = "320a0865-d16b-4ab5-9df0-33233ac5f952.md"
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("G3").Value = "320a0865-d16b-4ab5-9df0-33233ac5f952.af14c01cb35eab4cb4510f2c8ca8d566a39d3d5d.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-10-13 13:19:13"
$ws2.Range("J3").Value = "320a0865-d16b-4ab5-9df0-33233ac5f952.af14c01cb35eab4cb4510f2c8ca8d566a39d3d5d.zh-cn.xlf"
$ws2.Range("P3").Value = $errorDetail

$ws2.Range("A2").Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $urlTest0_320a, [Type]::Missing, [Type]::Missing, "524b8130-6d89-47c7-a944-7f48f88e90db.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("I2"), $urlZhcn_320a, [Type]::Missing, [Type]::Missing, "524b8130-6d89-47c7-a944-7f48f88e90db.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), $urlTest0_524b, [Type]::Missing, [Type]::Missing, "320a0865-d16b-4ab5-9df0-33233ac5f952.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("I3"), $urlZhcn_524b, [Type]::Missing, [Type]::Missing, "320a0865-d16b-4ab5-9df0-33233ac5f952.md") | Out-Null

$ws2.Columns.Item(16).ColumnWidth = 39.1

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "524b8130-6d89-47c7-a944-7f48f88e90db.md"
$ws3.Range("I2").Value = "524b8130-6d89-47c7-a944-7f48f88e90db.md"
$ws3.Range("G2").Value = "524b8130-6d89-47c7-a944-7f48f88e90db.29f77a1baee9a4f82de0e54920eab7c2a3320007.de-de.xlf"
$ws3.Range("J2").Value = "524b8130-6d89-47c7-a944-7f48f88e90db.29f77a1baee9a4f82de0e54920eab7c2a3320007.de-de.xlf"

$ws3.Range("A3").Value = "320a0865-d16b-4ab5-9df0-33233ac5f952.md"
$ws3.Range("I3").Value = "320a0865-d16b-4ab5-9df0-33233ac5f952.md"
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("G3").Value = "320a0865-d16b-4ab5-9df0-33233ac5f952.af14c01cb35eab4cb4510f2c8ca8d566a39d3d5d.de-de.xlf"
$ws3.Range("H3").Value = "2016-10-13 13:19:24"
$ws3.Range("J3").Value = "320a0865-d16b-4ab5-9df0-33233ac5f952.af14c01cb35eab4cb4510f2c8ca8d566a39d3d5d.de-de.xlf"
$ws3.Range("P3").Value = $errorDetail

$ws3.Range("A2").Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $urlTest0_320a, [Type]::Missing, [Type]::Missing, "524b8130-6d89-47c7-a944-7f48f88e90db.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("I2"), $urlDede_320a, [Type]::Missing, [Type]::Missing, "524b8130-6d89-47c7-a944-7f48f88e90db.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), $urlTest0_524b, [Type]::Missing, [Type]::Missing, "320a0865-d16b-4ab5-9df0-33233ac5f952.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("I3"), $urlDede_524b, [Type]::Missing, [Type]::Missing, "320a0865-d16b-4ab5-9df0-33233ac5f952.md") | Out-Null

$ws3.Columns.Item(16).ColumnWidth = 39.1
